# "Generate Report for handoff"
# b.md has been handed off again (new commit hash), so:
#  - Overview sheet: b.md row status -> "Ready for handoff" (both zh-cn/de-de columns)
#  - zh-cn detail sheet: status -> "Ready for handoff", new handoff file + datetime
#  - de-de detail sheet: status -> "Ready for handoff", new handoff file + datetime

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$zhHandoffFile = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$zhHandoffDate = "2016-02-16 09:48:24"
$deHandoffFile = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$deHandoffDate = "2016-02-16 09:48:37"

# --- Overview sheet: row 3 is b.md.md ---
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn detail sheet: row 3 is b.md.md ---
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = $zhHandoffFile
foreach ($hl in $wsZhCn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = $zhHandoffFile
    }
}
$wsZhCn.Range("D3").Value = $zhHandoffDate

# --- de-de detail sheet: row 3 is b.md.md ---
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("C3").Value = $deHandoffFile
foreach ($hl in $wsDeDe.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = $deHandoffFile
    }
}
$wsDeDe.Range("D3").Value = $deHandoffDate
